# Added JPMC UAT team demo script
# ---------------------------------------------------------------------------
# Adds four new test cases (JPMC-14, JPMC-15, JPMC-16, JPMC-17) to each of the
# three worksheets (default / en-us / es), fixes the jiraAssignee name on the
# JPMC-2 row, clears the stray jiraBugId placeholders on the en-us/es sheets,
# and widens the scenario/featureFile columns.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$assignee  = "Kiran Kumar V"
$email     = "kirankumar.baskar@photoninfotech.net"
$feature   = "JPMCSignin"
$featureCard = "Functional validation of Find a Card"

for ($i = 1; $i -le 3; $i++) {
    $ws = $wb.Worksheets.Item($i)

    # --- fix JPMC-2 row: jiraAssignee should read "Kiran Kumar V" ----------
    $ws.Range("D3").Value = $assignee

    # --- the en-us / es sheets had a leftover jiraBugId value on row 3 -----
    $ws.Range("C3").Value = ""

    # --- new row 4 : JPMC-14 -------------------------------------------------
    $ws.Range("A4").Value = "JPMC-14"
    $ws.Range("B4").Value = $email
    $ws.Range("D4").Value = $assignee
    $ws.Range("E4").Value = "Home Page Visual Validation"
    $ws.Range("F4").Value = $feature

    # --- new row 5 : JPMC-15 -------------------------------------------------
    $ws.Range("A5").Value = "JPMC-15"
    $ws.Range("B5").Value = $email
    $ws.Range("D5").Value = $assignee
    $ws.Range("E5").Value = "Page Load Performance Validation"
    $ws.Range("F5").Value = $feature

    # --- new row 6 : JPMC-16 -------------------------------------------------
    $ws.Range("A6").Value = "JPMC-16"
    $ws.Range("B6").Value = $email
    $ws.Range("D6").Value = $assignee
    $ws.Range("E6").Value = "Find a credit card for Personal Rewards Cash Back Balance Transfer"
    $ws.Range("F6").Value = $featureCard

    # --- new row 7 : JPMC-17 -------------------------------------------------
    $ws.Range("A7").Value = "JPMC-17"
    $ws.Range("B7").Value = $email
    $ws.Range("D7").Value = $assignee
    $ws.Range("E7").Value = "Find a credit card for Business Rewards Cash Back"
    $ws.Range("F7").Value = $featureCard

    # --- widen the scenario / featureFile columns ---------------------------
    $ws.Columns.Item(5).ColumnWidth = 58.0
    $ws.Columns.Item(6).ColumnWidth = 30.43
}

# The en-us sheet additionally keeps the JPMC-20 jira id on the new JPMC-16 row
$wsEnUs = $wb.Worksheets.Item(2)
$wsEnUs.Range("C6").Value = "JPMC-20"
